$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: force a run boundary at a given single point by toggling Bold
# on a 1-character probe range and reverting it. This produces a clean
# split with no residual direct formatting, because the final value is
# identical to the original (so Bold is omitted from the saved rPr).
# ---------------------------------------------------------------------
function Split-At($pos) {
    $probe = $d.Range($pos, $pos + 1)
    $orig = $probe.Font.Bold
    $probe.Font.Bold = $true
    $probe.Font.Bold = $orig
}

# =======================================================================
# Paragraph 1 : "Yerba Buena, 15 de Septiembre de 2011"
# =======================================================================
$p1 = $d.Paragraphs(1)
$p1.Format.KeepWithNext = $true
$p1.Format.SpaceAfter = 12

# =======================================================================
# Paragraph 2 : "ORDENANZA Nº 1805"
# =======================================================================
$p2 = $d.Paragraphs(2)
$p2.Format.KeepWithNext = $true
$p2.Format.SpaceBefore = 12
$p2.Format.SpaceAfter = 18
$p2.Range.Font.Bold = $true

# =======================================================================
# Paragraph 3 : "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA"
# Remove the leading space run, then format.
# =======================================================================
$p3 = $d.Paragraphs(3)
$lead = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
$lead.Text = ""

$p3 = $d.Paragraphs(3)
$p3.Format.KeepWithNext = $true
$p3.Format.SpaceBefore = 18
$p3.Format.SpaceAfter = 18
$p3.Format.LeftIndent = 99.2
$p3.Format.RightIndent = 99.2
$p3.Range.Font.Bold = $true

# =======================================================================
# Paragraph 4 : "ARTICULO PRIMERO: AUTORIZAR ..."
# Remove leading space; underline "ARTICULO PRIMERO" and ":" (as two
# separate runs); keepNext / spacing; drop justification; split the
# "Nº" occurrences into "N" / "º " / rest.
# =======================================================================
$p4 = $d.Paragraphs(4)
$lead = $d.Range($p4.Range.Start, $p4.Range.Start + 1)
$lead.Text = ""

$p4 = $d.Paragraphs(4)
$p4.Format.KeepWithNext = $true
$p4.Format.Alignment = 0
$p4.Format.SpaceAfter = 6

$base = $p4.Range.Start
$rTitle = $d.Range($base, $base + 17)              # "ARTICULO PRIMERO:"
Write-Output "p4 title=[$($rTitle.Text)]"
$rTitle.Font.Underline = 1
Split-At ($base + 16)                                # split "ARTICULO PRIMERO" | ":"

# --- split "Nº 26.445.871" -> "N" / "º " / "26.445.871..." -------------
$r = $d.Content
$r.Find.Execute("N", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $r.Find.Execute("D.N.I. N", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $r.End - 1
Write-Output "N-pos check=[$($d.Range($pos, $pos+2).Text)]"
Split-At ($pos)          # N | º 26...
Split-At ($pos + 2)      # º_ | 26...  (2 = 'º' + trailing space)

# --- split " Ubicado en Calle Pringle Nº 1263" --------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("Pringle N", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos2 = $r2.End - 1
Write-Output "pos2 check=[$($d.Range($pos2, $pos2+2).Text)]"
Split-At ($pos2)
Split-At ($pos2 + 2)

# --- split ", identificándose con el Padrón Nº 677.916" ----------------
$r3 = $d.Content
$null = $r3.Find.Execute("Padrón N", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos3 = $r3.End - 1
Write-Output "pos3 check=[$($d.Range($pos3, $pos3+2).Text)]"
Split-At ($pos3)
Split-At ($pos3 + 2)

# =======================================================================
# Paragraph 5 : "ARTICULO SEGUNDO: COMUNIQUESE, REGISTRESE Y ARCHIVESE."
# =======================================================================
$p5 = $d.Paragraphs(5)
$lead = $d.Range($p5.Range.Start, $p5.Range.Start + 1)
$lead.Text = ""

$p5 = $d.Paragraphs(5)
$p5.Format.KeepWithNext = $true
$p5.Format.Alignment = 0
$p5.Format.SpaceAfter = 6

$base5 = $p5.Range.Start
$rTitle5 = $d.Range($base5, $base5 + 17)            # "ARTICULO SEGUNDO:"
Write-Output "p5 title=[$($rTitle5.Text)]"
$rTitle5.Font.Underline = 1
Split-At ($base5 + 16)                               # split "ARTICULO SEGUNDO" | ":"

# =======================================================================
# Section properties: footer + starting page number
# =======================================================================
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$ftr.Range.InsertAfter("")
$ftr.PageNumbers.StartingNumber = 2481

Write-Output "ALL DONE"
